# Fruta / hortaliza, semanal
# Insert two new weekly price rows for "Palta" (Hass, Primera/Segunda, Peru
# origin, $/bandeja 10 kilos) dated 2022-09-22 (serial 44826) right above the
# existing row that used to be row 533, pushing all subsequent rows down by 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows at 533:534 (copies formatting, incl. the date number
# format on column D, from the row below).
$ws.Range("533:534").Insert()

# --- Row 533: Hass / Primera -------------------------------------------
$ws.Range("A533").Value = 7
$ws.Range("B533").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C533").Value = "Ñuble"
$ws.Range("D533").Value2 = 44826
$ws.Range("E533").Value = 16
$ws.Range("F533").Value = "Fruta"
$ws.Range("G533").Value = 100106
$ws.Range("H533").Value = "Oleaginosos"
$ws.Range("I533").Value = 100106002
$ws.Range("J533").Value = "Palta"
$ws.Range("K533").Value = "Hass"
$ws.Range("L533").Value = "Primera"
$ws.Range("M533").Value = 120
$ws.Range("N533").Value = 20000
$ws.Range("O533").Value = 23000
$ws.Range("P533").Value = 21500
$ws.Range("Q533").Value = "$/bandeja 10 kilos"
$ws.Range("R533").Value = "Perú"
$ws.Range("S533").Value = 2150
$ws.Range("T533").Value = 10

# --- Row 534: Hass / Segunda --------------------------------------------
$ws.Range("A534").Value = 7
$ws.Range("B534").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C534").Value = "Ñuble"
$ws.Range("D534").Value2 = 44826
$ws.Range("E534").Value = 16
$ws.Range("F534").Value = "Fruta"
$ws.Range("G534").Value = 100106
$ws.Range("H534").Value = "Oleaginosos"
$ws.Range("I534").Value = 100106002
$ws.Range("J534").Value = "Palta"
$ws.Range("K534").Value = "Hass"
$ws.Range("L534").Value = "Segunda"
$ws.Range("M534").Value = 120
$ws.Range("N534").Value = 20000
$ws.Range("O534").Value = 21000
$ws.Range("P534").Value = 20500
$ws.Range("Q534").Value = "$/bandeja 10 kilos"
$ws.Range("R534").Value = "Perú"
$ws.Range("S534").Value = 2050
$ws.Range("T534").Value = 10
